$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-18 21:05:57"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-18 21:05:52"
$wsZhCn.Range("K2").Value = "2016-08-18 21:06:17"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-18 21:05:57"
$wsDeDe.Range("K2").Value = "2016-08-18 21:06:24"
